$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "datos actualizados" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 06:20"

# --- Row 4: Estados Unidos stats update ---
$ws.Range("E4").Value = 155581
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 3166

# --- Rows 79/80: Uruguay <-> Kazajistan swap with updated stats ---
$ws.Range("A79").Value = "Kazajistan"
$ws.Range("B79").Value = 325
$ws.Range("C79").Value = 23
$ws.Range("D79").Value = 21
$ws.Range("E79").Value = 303
$ws.Range("F79").Value = 0

$ws.Range("A80").Value = "Uruguay"
$ws.Range("B80").Value = 320
$ws.Range("C80").Value = 0
$ws.Range("D80").Value = 25
$ws.Range("E80").Value = 294
$ws.Range("F80").Value = 9

# --- Rows 165-175: country reorder (Seychelles/Granada/Laos inserted earlier) ---
$ws.Range("A165").Value = "Seychelles"
$ws.Range("C165").Value = 2
$ws.Range("E165").Value = 10
$ws.Range("H165").Value = 0

$ws.Range("A166").Value = "Siria"
$ws.Range("D166").Value = 0
$ws.Range("H166").Value = 2

$ws.Range("A167").Value = "Groenlandia"
$ws.Range("B167").Value = 10
$ws.Range("D167").Value = 2
$ws.Range("E167").Value = 8

# Row 168 stays "Granada" (no change)

$ws.Range("A169").Value = "Laos"
$ws.Range("C169").Value = 1
$ws.Range("D169").Value = 0
$ws.Range("E169").Value = 9

$ws.Range("A170").Value = "Suazilandia"
$ws.Range("B170").Value = 9
$ws.Range("E170").Value = 9

$ws.Range("A171").Value = "Santa Lucia"
$ws.Range("B171").Value = 9
$ws.Range("D171").Value = 1

$ws.Range("A172").Value = "Libia"
$ws.Range("A173").Value = "Guinea-Bisau"
$ws.Range("A174").Value = "Surinam"
$ws.Range("A175").Value = "Mozambique"

# Row 176 stays "Guyana" (no change)
